$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "05/08/2025"
$ws.Range("A23").ClearFormats()
$ws.Range("B23").Value = "Alianza Huanuco"
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = "U. de Deportes"
$ws.Range("F23").Value = "L"
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 0.35
$ws.Range("L23").Value = 2.23
$ws.Range("M23").Value = 6
$ws.Range("N23").Value = 14
$ws.Range("O23").Value = 2
$ws.Range("P23").Value = 6
